{"js": "// Features-to-add list: remove two completed/obsolete items.\n//  - \"FIX CONTEXT MENU (...)\" is removed entirely. The \"_GoBack\" bookmark\n//    that lived inside that paragraph is preserved by re-inserting it at the\n//    start of the new first paragraph (\"Make upload routing in ajax.js\").\n//  - \"Make session_id global variable in DBManager\" is removed entirely.\n// All other bullet items keep their text/order unchanged.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst firstText = \"FIX CONTEXT MENU (SO THAT IT DOESN\\u2019T APPEAR ON BLACK SCREEN WHEN RIGHT MOUSE PRESSED)\";\nconst secondText = \"Make session_id global variable in DBManager\";\n\n// Does the paragraph carrying the \"_GoBack\" bookmark match the paragraph we\n// are about to delete? If so, remember it so we can re-home the bookmark.\nconst bookmarkRange = body.getBookmarkRangeOrNullObject(\"_GoBack\");\nbookmarkRange.load(\"isNullObject\");\nawait context.sync();\n\nlet bookmarkParaIndex = -1;\nif (!bookmarkRange.isNullObject) {\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    const loc = bookmarkRange.compareLocationWith(paragraphs.items[i]);\n    await context.sync();\n    if (loc.value === Word.LocationRelation.inside) {\n      bookmarkParaIndex = i;\n      break;\n    }\n  }\n}\n\n// Delete the two obsolete paragraphs (matched by their text so the script\n// does not depend on a fixed index).\nfor (const para of paragraphs.items) {\n  if (para.text === firstText || para.text === secondText) {\n    para.delete();\n  }\n}\nawait context.sync();\n\n// If the bookmark lived inside a paragraph we just removed, re-create it\n// (collapsed) at the very start of the new first paragraph.\nif (bookmarkParaIndex !== -1 && paragraphs.items[bookmarkParaIndex].text === firstText) {\n  const newParagraphs = body.paragraphs;\n  newParagraphs.load(\"items\");\n  await context.sync();\n  const startRange = newParagraphs.items[0].getRange(\"Start\");\n  startRange.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Features-to-add list: remove two completed/obsolete items.\n#  - \"FIX CONTEXT MENU (...)\" is removed entirely. The \"_GoBack\" bookmark that\n#    lived at the end of that paragraph's text is preserved by relocating it\n#    to the start of the new first paragraph (\"Make upload routing in ajax.js\").\n#  - \"Make session_id global variable in DBManager\" is removed entirely.\n# All other bullet items keep their text/order unchanged.\n\n$d = $word.ActiveDocument\n\n# --- Remove \"FIX CONTEXT MENU ...\" while preserving the _GoBack bookmark ---\n$first = $d.Paragraphs.Item(1).Range\n$firstText = $first.Text\n# Trim the trailing paragraph mark so we only delete the visible text and\n# leave the (now empty) paragraph + its collapsed bookmark in place.\n$textLen = $firstText.Length - 1\n$d.Range($first.Start, $first.Start + $textLen).Delete()\n\n# Re-insert the replacement text after the now-empty paragraph 1 so the\n# bookmark (still collapsed at the paragraph start) ends up ahead of it.\n$first = $d.Paragraphs.Item(1).Range\n$first.InsertAfter(\"Make upload routing in ajax.js\")\n$first = $d.Paragraphs.Item(1).Range\n$first.LanguageID = \"en-US\"\n\n# The old \"Make upload routing in ajax.js\" paragraph (now paragraph 2, a\n# duplicate of paragraph 1) is removed.\n$d.Paragraphs.Item(2).Range.Delete()\n\n# --- Remove \"Make session_id global variable in DBManager\" ---\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq \"Make session_id global variable in DBManager\") {\n        $p.Range.Delete()\n        break\n    }\n}\n"}
